# aci-tenant-deploy input spreadsheet: add a "Create Node Profiles" block to
# the L3Out sheet, retarget the L3Out header row to L3Out VRF / L3Out L3
# Domain / Description, fill the Routing Protocol value, and make L3Out the
# active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L3Out")

# --- Row 1 title -----------------------------------------------------
$ws.Range("A1").Value = "Create L3Outs"

# --- Row 3 column headers (L3Out VRF / L3Out L3 Domain / Description) -
$ws.Range("D3").Value = "L3Out VRF"
$ws.Range("E3").Value = "L3Out L3 Domain"
$ws.Range("F3").Value = "Description"

# --- Row 4 sample data: fill in the Routing Protocol (now "L3Out L3 Domain" col E) value
$ws.Range("E4").Value = "bgp"

# --- Build the new "Create Node Profiles" block (rows 13-23) by copying
# the formatting of the existing "Create L3Outs" block (rows 1-11) ------
$ws.Range("A1:M11").Copy()
$ws.Range("A13:M23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(13).RowHeight = 20.25
$ws.Rows.Item(14).RowHeight = 35.25
$ws.Rows.Item(14).Rows.AutoFit()
$ws.Rows.Item(15).RowHeight = 17.25
for ($r = 16; $r -le 23; $r++) {
    $ws.Rows.Item($r).RowHeight = 16.5
}

$ws.Range("A13:M13").Merge()
$ws.Range("B14:M14").Merge()

# Section header + instructional text
$ws.Range("A13").Value = "Create Node Profiles"
$ws.Range("B14").Value = "Type tnt_add will only Create a Tenant; `nType tnt_vrf will create the Tenant and a vrf with the name {tenant}_vrf in the common Tenant.  We recommend using this in most cases"

# Column headers (row 15)
$ws.Range("A15").Value = "Type"
$ws.Range("B15").Value = "Leaf Name"
$ws.Range("C15").Value = "Leaf Node Id"
$ws.Range("D15").Value = "Interface Type"
$ws.Range("E15").Value = "Interface"
$ws.Range("F15").Value = "VLAN"
$ws.Range("G15").Value = "IPv4/Prefix"

# Sample data row (row 16)
$ws.Range("A16").Value = "nodeprof"
$ws.Range("B16").Value = "leaf201"
$ws.Range("C16").Value = "201"
$ws.Range("D16").Value = "Routed Interface"
$ws.Range("E16").Value = "1/1"

# --- Data validation ---------------------------------------------------
# Remove the old "vpc_pair" list validation on A5:A11 (the VPC workflow no
# longer applies to this sheet).
$ws.Range("A5:A11").Validation.Delete()

# Extend the Type-column list validation down across the whole L3Out block.
$ws.Range("A4").Validation.Delete()
$ws.Range("A4:A11").Validation.Add(3, 1, 1, """l3out""")

# Drop the old bgp/ospf validation on E4 (E4 now just holds a plain value).
$ws.Range("E4").Validation.Delete()

# VPC-ID style whole-number validation also applies to the new block's blank rows.
$ws.Range("B17").Validation.Add(1, 1, 1, 1, 1000)
$ws.Range("B17").Validation.ErrorTitle = "VPC ID"
$ws.Range("B17").Validation.ErrorMessage = "The VPC ID must be a number between 1 and 1000"
$ws.Range("B19").Validation.Add(1, 1, 1, 1, 1000)
$ws.Range("B19").Validation.ErrorTitle = "VPC ID"
$ws.Range("B19").Validation.ErrorMessage = "The VPC ID must be a number between 1 and 1000"
$ws.Range("B21:B23").Validation.Add(1, 1, 1, 1, 1000)
$ws.Range("B21:B23").Validation.ErrorTitle = "VPC ID"
$ws.Range("B21:B23").Validation.ErrorMessage = "The VPC ID must be a number between 1 and 1000"

# New list validations for the Node Profile block.
$ws.Range("A16:A23").Validation.Add(3, 1, 1, """nodeprof""")
$ws.Range("D16").Validation.Add(3, 1, 1, """Routed Interface,sub-interface,svi""")

# --- Page setup ----------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Make L3Out the active / selected sheet -----------------------------
$ws.Activate()
$ws.Range("F16").Select()

Write-Output "edit complete"
